# Fixed errors in moments computation.
#
# 1) "data" sheet: update simulated-moment values (columns D/E/F, rows 5-23)
#    to the corrected figures.
# 2) "table" sheet: the % rows (11,12,13,14,19,20) were double-multiplying
#    by 100 in the "Data" column - drop the stray "*100".
# 3) "table" sheet: a new (currently blank) formatted row 24 / cell F24
#    appears, extending the used range.
# 4) View-state touch-up: "data" tab ends up the active/selected tab with
#    its selection on C20:J20; "table" tab is no longer the selected tab
#    and its selection collapses to F22.

$wb = $excel.ActiveWorkbook

$data = $wb.Worksheets.Item("data")
$table = $wb.Worksheets.Item("table")

# --- 1) corrected simulation moments on "data" ---------------------------

$data.Range("D5").Value = 2.486743561837384
$data.Range("E5").Value = 2.473010797023773
$data.Range("F5").Value = 0.07832692683838968

$data.Range("D6").Value = 0.07191216137102167
$data.Range("E6").Value = 0.0737045496404171
$data.Range("F6").Value = 0.002791894768823342

$data.Range("D7").Value = -0.1763006811907082
$data.Range("E7").Value = -0.3557592459619045
$data.Range("F7").Value = 0.01269640106306268

$data.Range("D8").Value = 0.3869429884421218
$data.Range("E8").Value = 0.2193040450364351
$data.Range("F8").Value = 0.008640766787331089

$data.Range("D9").Value = 2.500736553924154
$data.Range("E9").Value = 2.58952593588829
$data.Range("F9").Value = 0.08228883058330325

$data.Range("D10").Value = 0.2858268679925717
$data.Range("E10").Value = 0.3511869743466377
$data.Range("F10").Value = 0.01197301821384926

$data.Range("D11").Value = 2.464477557672385
$data.Range("E11").Value = 2.492690888904214
$data.Range("F11").Value = 0.07911715410396938

$data.Range("D12").Value = 0.01881308033911991
$data.Range("E12").Value = 0.0388672060009751
$data.Range("F12").Value = 0.004052441122115701

$data.Range("D13").Value = 0.3297214372224465
$data.Range("E13").Value = 0.282855486013866
$data.Range("F13").Value = 0.01284582756155296

$data.Range("D14").Value = 0.3625353249899071
$data.Range("E14").Value = 0.5285845146094897
$data.Range("F14").Value = 0.01950110813856742

$data.Range("D15").Value = 0.2689301574485264
$data.Range("E15").Value = 0.1486927933756691
$data.Range("F15").Value = 0.008571077048086666

$data.Range("D16").Value = 0.2674542752462063
$data.Range("E16").Value = 0.2381824513994524
$data.Range("F16").Value = 0.0193161450629829

$data.Range("D17").Value = 0.1763744298772759
$data.Range("E17").Value = 0.2242545439751264
$data.Range("F17").Value = 0.02157942660681845

$data.Range("D18").Value = -0.0125117225943009
$data.Range("E18").Value = -0.06574915627640092
$data.Range("F18").Value = 0.02065201320820884

$data.Range("D19").Value = -0.003802775894403772
$data.Range("E19").Value = -0.02422105296519655
$data.Range("F19").Value = 0.0203809442777414

$data.Range("E20").Value = 0.3770313102118873
$data.Range("F20").Value = 0.01442944072430714

$data.Range("D21").Value = 0
$data.Range("E21").Value = 0.585378801634811
$data.Range("F21").Value = 0.02020568064632542

$data.Range("D22").Value = 0.4657037156153698
$data.Range("E22").Value = 0.4184576788780685
$data.Range("F22").Value = 0.01973609522177048

$data.Range("D23").Value = 0.1967505437285323
$data.Range("E23").Value = 0.2073579741848141
$data.Range("F23").Value = 0.01427416787326762

# --- 2) drop the stray "*100" on the "table" sheet's Data column ---------

$table.Range("F11").Formula = "=data!E12"
$table.Range("F12").Formula = "=data!E13"
$table.Range("F13").Formula = "=data!E14"
$table.Range("F14").Formula = "=data!E15"
$table.Range("F19").Formula = "=data!E20"
$table.Range("F20").Formula = "=data!E21"

# --- 3) new blank, formatted row 24 on "table" ----------------------------
# (mirrors the "0.00"-style, centred formatting used by the rest of column F)

$table.Range("F24").NumberFormat = "0.00"
$table.Range("F24").HorizontalAlignment = -4108

# --- 4) view-state: "table" selection collapses, "data" becomes active ---

$table.Range("F22").Select()

$data.Activate()
$data.Range("C20:J20").Select()
